$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "[Valmir-Calderaria-3B, -, -, -]"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("E3").Value = "['MEC-3B-Usin. CNC', 'MEC-2B-Ajustagem', -, -]"
$ws.Range("F3").Value = "-"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("E4").Value = "['MEC-3B-Usin. CNC', 'MEC-2B-Ajustagem', -, -]"
$ws.Range("F4").Value = "-"

# Row 6
$ws.Range("B6").Value = "-"
$ws.Range("E6").Value = "[-, 'MEC-2B-Ajustagem', -, -]"
$ws.Range("F6").Value = "-"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("E7").Value = "[-, 'MEC-2B-Ajustagem', -, -]"

# Row 8
$ws.Range("E8").Value = "[-, -, -, 'MEC-3B-Usin. CNC']"

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("D18").Value = "[-, -, -, 'MEC-2NB-Usin. CNC']"
$ws.Range("E18").Value = "[-, 'MEC-2NB-Usin. CNC', 'MEC-1NA-Ajustagem', 'MEC-1NB-Ajustagem']"
$ws.Range("F18").Value = "-"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("D19").Value = "[-, -, -, 'MEC-2NB-Usin. CNC']"
$ws.Range("E19").Value = "[-, -, 'MEC-1NA-Ajustagem', 'MEC-1NB-Ajustagem']"

# Row 20
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "[-, -, -, 'MEC-2NB-Usin. CNC']"
$ws.Range("E20").Value = "[-, -, 'MEC-1NA-Ajustagem', 'MEC-1NB-Ajustagem']"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "[-, -, 'MEC-1NA-Ajustagem', 'MEC-1NB-Ajustagem']"
